$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Pneumonia / Cold columns)
$ws.Range("F2").Value = 0.01
$ws.Range("H2").Value = 0.99

# Row 3
$ws.Range("F3").Value = 0.01
$ws.Range("H3").Value = 0.01

# Row 4
$ws.Range("F4").Value = 0.99
$ws.Range("G4").Value = 0.01
$ws.Range("H4").Value = 0.01

# Row 5
$ws.Range("F5").Value = 0.01

# Row 6
$ws.Range("F6").Value = 0.01
$ws.Range("H6").Value = 0.99
